$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3497
$ws.Range("F5").Value = 8199
$ws.Range("F7").Value = 81
$ws.Range("F8").Value = 2153
$ws.Range("F10").Value = 181
$ws.Range("F12").Value = 1184
$ws.Range("F15").Value = 510
$ws.Range("F16").Value = 566
$ws.Range("F17").Value = 71
$ws.Range("F18").Value = 48
$ws.Range("F19").Value = 432
$ws.Range("F21").Value = 7073
$ws.Range("F22").Value = 143
$ws.Range("F23").Value = 54770
$ws.Range("F24").Value = 54770
$ws.Range("F25").Value = 4338
$ws.Range("F27").Value = 843
$ws.Range("F28").Value = 406
$ws.Range("F29").Value = 87
$ws.Range("F30").Value = 872
$ws.Range("F32").Value = 593
$ws.Range("F33").Value = 2482
$ws.Range("F34").Value = 577
$ws.Range("F35").Value = 20
$ws.Range("F36").Value = 18
$ws.Range("F37").Value = 856
$ws.Range("F38").Value = 1151
$ws.Range("F39").Value = 768
$ws.Range("F40").Value = 143
$ws.Range("F42").Value = 1060
$ws.Range("F44").Value = 758
$ws.Range("F45").Value = 144
$ws.Range("F47").Value = 139
$ws.Range("F48").Value = 2

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 54
$ws.Range("F8").Value = 164
$ws.Range("F9").Value = 28
$ws.Range("F10").Value = 47
$ws.Range("F12").Value = 103
$ws.Range("F14").Value = 38
$ws.Range("F16").Value = 7439
$ws.Range("F17").Value = 102
$ws.Range("F32").Value = 19
$ws.Range("F39").Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2269
$ws.Range("F5").Value = 1532
$ws.Range("F7").Value = 649
$ws.Range("F9").Value = 9316
$ws.Range("F10").Value = 1625
$ws.Range("F12").Value = 74
$ws.Range("F15").Value = 143

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3497
$ws.Range("F4").Value = 8199
$ws.Range("F5").Value = 649
$ws.Range("F7").Value = 74
$ws.Range("F9").Value = 81
$ws.Range("F12").Value = 181
$ws.Range("F13").Value = 1184
$ws.Range("F15").Value = 143
$ws.Range("F17").Value = 71
$ws.Range("F18").Value = 432
$ws.Range("F19").Value = 7073
$ws.Range("F20").Value = 143
$ws.Range("F21").Value = 54770
$ws.Range("F22").Value = 164
$ws.Range("F23").Value = 164
$ws.Range("F24").Value = 28
$ws.Range("F25").Value = 47
$ws.Range("F26").Value = 4338
$ws.Range("F29").Value = 843
$ws.Range("F30").Value = 406
$ws.Range("F31").Value = 593
$ws.Range("F32").Value = 103
$ws.Range("F34").Value = 2482
$ws.Range("F35").Value = 577
$ws.Range("F36").Value = 38
$ws.Range("F37").Value = 18
$ws.Range("F38").Value = 856
$ws.Range("F39").Value = 1151
$ws.Range("F40").Value = 102
$ws.Range("F42").Value = 143
$ws.Range("F43").Value = 1060
$ws.Range("F46").Value = 758
$ws.Range("F47").Value = 144
$ws.Range("F49").Value = 140
$ws.Range("F51").Value = 21
